$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.149.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = "'2.244.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'306.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").Value = "'95.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").Value = "'0.573"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = "'0.525"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = "'34.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.64%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").Value = "'7.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").Value = "'2.585.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").Value = "'2.254.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").Value = "'0.833"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = "'13.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").Value = "'44.038.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").Value = "'0.0₃0975"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("D20").Value = "'6.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = "'12.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.44%  '
$ws.Range("D22").Value = "'65.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = "'236.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = "'37.48"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.09%  '
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("D32").Value = "'151.74"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.57%  '
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("E34").Value = '  +4.10%  '
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("D36").Value = "'0.109"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("E38").Value = '  -6.28%  '
$ws.Range("E39").Value = '  -5.17%  '
$ws.Range("E40").Value = '  -3.87%  '
$ws.Range("D41").Value = "'14.43"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -8.84%  '
$ws.Range("D42").Value = "'0.0298"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = "'1.741.58"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = "'83.12"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.85%  '
$ws.Range("D46").Value = "'0.190"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").Value = "'100.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Value = "'8.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("D50").Value = "'54.59"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'67.79"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.21%  '
